# Add 2022-Q3 data:
#  - shift the "总计" (summary) sheet rows down by one and insert a new
#    2022-Q3 row at the top of the data (row 2)
#  - insert a brand-new "2022-Q3" worksheet (cloned from the "2021-Q4"
#    sheet so it keeps the same layout/formatting) right before the
#    "2021-Q4" tab, then fill it with the 2022-Q3 holdings data

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" summary sheet: insert a new first data row for
#    2022-Q3 and push the existing rows (2021-Q4, 2021-Q3, 2021-Q2,
#    2020-Q4) down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Grab the old values (rows 2-5) before we overwrite anything.
# NB: read through .Value2 (not .Value) - this COM shim's .Value getter
# does not round-trip through a PowerShell variable correctly.
$oldB2 = $summary.Cells.Item(2,2).Value2
$oldC2 = $summary.Cells.Item(2,3).Value2
$oldD2 = $summary.Cells.Item(2,4).Value2
$oldB3 = $summary.Cells.Item(3,2).Value2
$oldC3 = $summary.Cells.Item(3,3).Value2
$oldD3 = $summary.Cells.Item(3,4).Value2
$oldB4 = $summary.Cells.Item(4,2).Value2
$oldC4 = $summary.Cells.Item(4,3).Value2
$oldD4 = $summary.Cells.Item(4,4).Value2
$oldB5 = $summary.Cells.Item(5,2).Value2
$oldC5 = $summary.Cells.Item(5,3).Value2
$oldD5 = $summary.Cells.Item(5,4).Value2

# Copy the style of the existing row 5 down into the brand new row 6
# (column A carries the bordered/bold "index" style).
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)

$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(6,2).Value = $oldB5
$summary.Cells.Item(6,3).Value = $oldC5
$summary.Cells.Item(6,4).Value = $oldD5

$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = $oldB4
$summary.Cells.Item(5,3).Value = $oldC4
$summary.Cells.Item(5,4).Value = $oldD4

$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = $oldB3
$summary.Cells.Item(4,3).Value = $oldC3
$summary.Cells.Item(4,4).Value = $oldD3

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = $oldB2
$summary.Cells.Item(3,3).Value = $oldC2
$summary.Cells.Item(3,4).Value = $oldD2

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 11
$summary.Cells.Item(2,4).Value = 2.04

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet before "2021-Q4", cloning the
#    "2021-Q4" sheet so formatting (borders/bold header/index column)
#    matches the other quarterly sheets.
# ---------------------------------------------------------------------
$q4sheet = $wb.Worksheets.Item("2021-Q4")
$q4sheet.Copy($q4sheet)

# Re-fetch "2021-Q4" by name: the original sheet's Index is stale right
# after Copy(), but Item("2021-Q4") still resolves to the original sheet
# (now shifted one slot to the right), so the clone sits right before it.
$q4sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Item($q4sheet.Index - 1)
$newSheet.Name = "2022-Q3"

# Extend column-A style (bordered/bold "index" style) and header style
# (bold/centered/bordered) down/across to cover the extra rows/columns
# this sheet needs (12 rows vs the cloned sheet's 5).
$newSheet.Range("A2:A5").Copy()
$newSheet.Range("A6:A12").PasteSpecial(-4122)
$newSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Header row
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Columns B,D,E,F,G hold plain text (e.g. "008188", "21.72"), not real
# numbers, so force text formatting before assigning them (otherwise
# Excel auto-coerces numeric-looking strings into numbers).
$newSheet.Range("B2:B12").NumberFormat = "@"
$newSheet.Range("D2:G12").NumberFormat = "@"

$rows = @(
    @(0, "008188", "前海开源稳健增长三年持有期混合", "21.72", "66.43", "2.87", "0.6234", 7),
    @(1, "002258", "大成国企改革灵活配置混合",       "17.20", "90.32", "2.92", "0.5022", 10),
    @(2, "010826", "大成产业趋势混合A",               "11.59", "90.99", "3.82", "0.4427", 8),
    @(3, "010827", "大成产业趋势混合C",               "6.29",  "90.99", "3.82", "0.2403", 8),
    @(4, "011287", "前海开源聚慧三年持有期混合",       "2.84",  "66.61", "2.89", "0.0821", 7),
    @(5, "006775", "前海开源优质成长混合",             "2.48",  "69.80", "2.94", "0.0729", 8),
    @(6, "006216", "前海开源价值成长灵活配置混合A",     "1.11",  "65.31", "2.77", "0.0307", 8),
    @(7, "002407", "前海开源恒远灵活配置混合",         "1.03",  "67.57", "2.90", "0.0299", 7),
    @(8, "006217", "前海开源价值成长灵活配置混合C",     "0.47",  "65.31", "2.77", "0.0130", 8),
    @(9, "002020", "国都创新驱动灵活配置混合",         "0.12",  "65.45", "4.77", "0.0057", 2),
    @(10, "005247", "国都量化精选混合",                "0.02",  "75.50", "2.32", "0.0005", 9)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $newSheet.Cells.Item($r,1).Value = $data[0]
    $newSheet.Cells.Item($r,2).Value = $data[1]
    $newSheet.Cells.Item($r,3).Value = $data[2]
    $newSheet.Cells.Item($r,4).Value = $data[3]
    $newSheet.Cells.Item($r,5).Value = $data[4]
    $newSheet.Cells.Item($r,6).Value = $data[5]
    $newSheet.Cells.Item($r,7).Value = $data[6]
    $newSheet.Cells.Item($r,8).Value = $data[7]
}

# The "text-format" trick above leaves the cells tagged with an explicit
# (non-bordered) text style; reset B2:B12/D2:G12 back to the workbook's
# default (un-styled) look by pasting the format from a genuinely blank,
# default-styled cell over them - this keeps the *values* as text while
# dropping the incidental style index.
$newSheet.Cells.Item(1,1).Copy()
$newSheet.Range("B2:B12").PasteSpecial(-4122)
$newSheet.Range("D2:G12").PasteSpecial(-4122)

# Restore the originally-active tab ("2020-Q4", the last sheet) so the
# new/duplicated sheets we touched aren't left marked as the selected tab.
$wb.Worksheets.Item("2020-Q4").Activate()

